{"js": "// Replace the three-digit-by-one-digit multiplication prompts in the\n// worksheet table with a new set of problems, matching the author's\n// edit (old expression text -> new expression text).\nconst replacements = [\n  [\"204\u00d73=\", \"255\u00d74=\"],\n  [\"509\u00d75=\", \"267\u00d74=\"],\n  [\"268\u00d74=\", \"217\u00d76=\"],\n  [\"537\u00d74=\", \"954\u00d79=\"],\n  [\"993\u00d72=\", \"529\u00d76=\"],\n  [\"492\u00d78=\", \"261\u00d75=\"],\n  [\"767\u00d75=\", \"672\u00d73=\"],\n  [\"424\u00d74=\", \"113\u00d73=\"],\n  [\"797\u00d76=\", \"318\u00d78=\"],\n  [\"441\u00d75=\", \"222\u00d76=\"],\n  [\"607\u00d72=\", \"345\u00d74=\"],\n  [\"265\u00d73=\", \"101\u00d75=\"],\n  [\"720\u00d74=\", \"198\u00d75=\"],\n  [\"473\u00d78=\", \"846\u00d77=\"],\n  [\"299\u00d74=\", \"723\u00d78=\"],\n  [\"843\u00d74=\", \"580\u00d76=\"],\n  [\"148\u00d73=\", \"484\u00d73=\"],\n  [\"194\u00d77=\", \"394\u00d78=\"],\n  [\"393\u00d78=\", \"177\u00d79=\"],\n  [\"492\u00d73=\", \"609\u00d75=\"],\n  [\"930\u00d77=\", \"202\u00d74=\"],\n  [\"692\u00d79=\", \"206\u00d75=\"],\n  [\"348\u00d76=\", \"212\u00d76=\"],\n  [\"388\u00d79=\", \"231\u00d77=\"],\n  [\"668\u00d74=\", \"204\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication prompts in the\n# worksheet table with a new set of problems, matching the author's\n# edit (old expression text -> new expression text).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"204\u00d73=\", \"255\u00d74=\"),\n    @(\"509\u00d75=\", \"267\u00d74=\"),\n    @(\"268\u00d74=\", \"217\u00d76=\"),\n    @(\"537\u00d74=\", \"954\u00d79=\"),\n    @(\"993\u00d72=\", \"529\u00d76=\"),\n    @(\"492\u00d78=\", \"261\u00d75=\"),\n    @(\"767\u00d75=\", \"672\u00d73=\"),\n    @(\"424\u00d74=\", \"113\u00d73=\"),\n    @(\"797\u00d76=\", \"318\u00d78=\"),\n    @(\"441\u00d75=\", \"222\u00d76=\"),\n    @(\"607\u00d72=\", \"345\u00d74=\"),\n    @(\"265\u00d73=\", \"101\u00d75=\"),\n    @(\"720\u00d74=\", \"198\u00d75=\"),\n    @(\"473\u00d78=\", \"846\u00d77=\"),\n    @(\"299\u00d74=\", \"723\u00d78=\"),\n    @(\"843\u00d74=\", \"580\u00d76=\"),\n    @(\"148\u00d73=\", \"484\u00d73=\"),\n    @(\"194\u00d77=\", \"394\u00d78=\"),\n    @(\"393\u00d78=\", \"177\u00d79=\"),\n    @(\"492\u00d73=\", \"609\u00d75=\"),\n    @(\"930\u00d77=\", \"202\u00d74=\"),\n    @(\"692\u00d79=\", \"206\u00d75=\"),\n    @(\"348\u00d76=\", \"212\u00d76=\"),\n    @(\"388\u00d79=\", \"231\u00d77=\"),\n    @(\"668\u00d74=\", \"204\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
